$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new 2-row block above the old row 11 (pushes indent_1..indent_5
# rows down from 11/13/15/17/19 to 13/15/17/19/21, leaving the new blank
# row 12 as the spacer that already separated every other entry).
$ws.Rows("11:12").Insert()

# Populate the new header-style row (indent_0) that now lives at row 11.
$ws.Range("A11").Value = "indent_0"
$ws.Range("B11").Value = 19

# Style A11 like the "title" cells (16pt font) on a white themed fill -
# this allocates the new font/fill combination and the new cellXfs entry.
$ws.Range("A11").Font.Size = 16
$ws.Range("A11").Interior.Pattern = 1
$ws.Range("A11").Interior.ThemeColor = 2
$ws.Range("A11").Interior.TintAndShade = 0

# Row 11 is a bit taller than the other indent rows.
$ws.Rows("11").RowHeight = 20

# The shifted-down rows keep their original styles/labels (handled by the
# row insert above) but their B values step down by 2 instead of by 1.
$ws.Range("B13").Value = 17
$ws.Range("B15").Value = 15
$ws.Range("B17").Value = 13
$ws.Range("B19").Value = 11
$ws.Range("B21").Value = 9

# Match the saved selection from the authored workbook.
$ws.Range("B22").Select()
